$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test case description
$ws.Range("B1").Value = "Test Case: Testing to see any errors while logging in"

# Step 1
$ws.Range("C2").Value = "Step 1: Sign-in with no fields filled with information"
$ws.Range("D2").Value = "An error text pops up saying that I need to fill out those field"

# Step 2
$ws.Range("C3").Value = "Step 2: Sign in with email filled out but no password"
$ws.Range("D3").Value = "An error text pops up saying that I need to fill out the password field"

# Step 3
$ws.Range("C4").Value = "Step 3: Sign in with password filled out but not email"
$ws.Range("D4").Value = "An error text pops up saying that the email is required"

# Step 4
$ws.Range("C5").Value = "Step 4: Fill out both fields with the wrong information"
$ws.Range("D5").Value = "An error pops up saying that the email or password is wrong"

# Step 5
$ws.Range("C6").Value = "Step 5: Fill out the fields with the correct information "
$ws.Range("D6").Value = "I am logged in and redirected to the user's dashboard with no errors."

# Update the selection shown in the sheet view
$ws.Range("B1:F6").Select()
